# Translate the "fdr" (French, facteur de risque) entries in column E
# (the "table2" column) to the English "risk_factor", and move the
# active-cell selection to E15 (matching the saved selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("E2", "E3", "E4", "E6", "E9", "E10", "E11", "E12", "E15")
foreach ($addr in $cells) {
    $ws.Range($addr).Value = "risk_factor"
}

[void]$ws.Range("E15").Select()
